$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.591.17"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.538.47"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'313.34"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "'95.05"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("D7").Value = "'0.578"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "'36.26"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'0.0815"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "'7.71"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "2.928.71"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "'15.68"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "2.517.10"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "'0.866"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "42.658.26"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'13.05"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").Value = "'6.68"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "0.0₃0970"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").Value = "'71.13"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").Value = "'254.96"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "'27.42"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "'0.993"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +11.76%  "
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("D30").Value = "'10.06"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").Value = "'155.62"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("D33").Value = "'19.82"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").Value = "'3.41"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "'0.0795"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "'2.63"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("D39").Value = "'24.77"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("D40").Value = "'0.119"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'2.21"
$ws.Range("E41").Value = "  +6.24%  "
$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "'3.84"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "'0.0303"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "2.040.33"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").Value = "'85.20"
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("D48").Value = "'8.93"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'75.58"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("D50").Value = "2.784.19"
$ws.Range("E51").Value = "  -0.05%  "
